{"js": "// Office.js (Word JavaScript API) script applying the r\u00e9sum\u00e9 update:\n//  - Name line: \"DHEERAJ CHAND\" -> \"Dheeraj Chand\"\n//  - Title line: \"Data Analytics & Technology Consultant\" -> \"Professional Title\"\n//  - Contact line: \"(202) 550-7110 | Dheeraj.Chand@gmail.com\" -> \"202.550.7110 | dheeraj.chand@gmail.com\"\n//  - Summary: \"20+ years\" -> \"21 years\"\n//  - Employer line: \"Siege Analytics, Austin, TX | 2005 \u2013 Present\" -> \"Your Company Name, Your City, ST | 2005 \u2013 Present\"\n//  - Remove the \"DATA PRODUCTS MANAGER\" and \"TECHNICAL CONSULTANT\" job blocks entirely\n//  - Remove the \"Client Impact\" and \"Methodological Innovation\" blocks entirely\n\nconst body = context.document.body;\n\nasync function replaceFirst(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replacement, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 1. Simple text replacements (formatting of the surrounding run is preserved).\nawait replaceFirst(\"DHEERAJ CHAND\", \"Dheeraj Chand\");\nawait replaceFirst(\"Data Analytics & Technology Consultant\", \"Professional Title\");\nawait replaceFirst(\"(202) 550-7110 | Dheeraj.Chand@gmail.com\", \"202.550.7110 | dheeraj.chand@gmail.com\");\nawait replaceFirst(\"20+ years\", \"21 years\");\nawait replaceFirst(\"Siege Analytics, Austin, TX | 2005\", \"Your Company Name, Your City, ST | 2005\");\n\n// 2. Remove whole paragraph blocks.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfunction paraIndex(matchText) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.trim() === matchText) {\n      return i;\n    }\n  }\n  return -1;\n}\n\nfunction collectRange(startText, endText) {\n  const startIdx = paraIndex(startText);\n  const endIdx = paraIndex(endText);\n  const collected = [];\n  if (startIdx !== -1 && endIdx !== -1 && startIdx <= endIdx) {\n    for (let i = startIdx; i <= endIdx; i++) {\n      collected.push(paragraphs.items[i]);\n    }\n  }\n  return collected;\n}\n\n// Block 1: \"DATA PRODUCTS MANAGER\" ... \"\u25b8 Implemented scalable telephony integration systems for large-scale data collection\"\nconst block1 = collectRange(\n  \"DATA PRODUCTS MANAGER\",\n  \"\u25b8 Implemented scalable telephony integration systems for large-scale data collection\"\n);\n\n// Block 2: \"Client Impact\" ... \"\u2713 Established best practices for multi-tenant data architecture and security compliance\"\nconst block2 = collectRange(\n  \"Client Impact\",\n  \"\u2713 Established best practices for multi-tenant data architecture and security compliance\"\n);\n\nfor (const p of block1) {\n  p.delete();\n}\nfor (const p of block2) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "# Word COM interop script applying the r\u00e9sum\u00e9 update:\n#  - Name line: \"DHEERAJ CHAND\" -> \"Dheeraj Chand\"\n#  - Title line: \"Data Analytics & Technology Consultant\" -> \"Professional Title\"\n#  - Contact line: \"(202) 550-7110 | Dheeraj.Chand@gmail.com\" -> \"202.550.7110 | dheeraj.chand@gmail.com\"\n#  - Summary: \"20+ years\" -> \"21 years\"\n#  - Employer line: \"Siege Analytics, Austin, TX | 2005 - Present\" -> \"Your Company Name, Your City, ST | 2005 - Present\"\n#  - Remove the \"DATA PRODUCTS MANAGER\" and \"TECHNICAL CONSULTANT\" job blocks entirely\n#  - Remove the \"Client Impact\" and \"Methodological Innovation\" blocks entirely\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-Text \"DHEERAJ CHAND\" \"Dheeraj Chand\"\nReplace-Text \"Data Analytics & Technology Consultant\" \"Professional Title\"\nReplace-Text \"(202) 550-7110 | Dheeraj.Chand@gmail.com\" \"202.550.7110 | dheeraj.chand@gmail.com\"\nReplace-Text \"20+ years\" \"21 years\"\nReplace-Text \"Siege Analytics, Austin, TX | 2005\" \"Your Company Name, Your City, ST | 2005\"\n\n# Remove the \"DATA PRODUCTS MANAGER\" .. \"Implemented scalable telephony integration systems for large-scale data collection\" block\n$startPara = $null\n$endPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"DATA PRODUCTS MANAGER\") { $startPara = $p }\n    if ($t -eq \"\u25b8 Implemented scalable telephony integration systems for large-scale data collection\") { $endPara = $p }\n}\nif ($startPara -ne $null -and $endPara -ne $null) {\n    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $r.Delete()\n}\n\n# Remove the \"Client Impact\" .. \"Established best practices for multi-tenant data architecture and security compliance\" block\n$startPara2 = $null\n$endPara2 = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Client Impact\") { $startPara2 = $p }\n    if ($t -eq \"\u2713 Established best practices for multi-tenant data architecture and security compliance\") { $endPara2 = $p }\n}\nif ($startPara2 -ne $null -and $endPara2 -ne $null) {\n    $r2 = $d.Range($startPara2.Range.Start, $endPara2.Range.End)\n    $r2.Delete()\n}\n"}
